# TestCoverage.xlsx - update method names / load test data
# (see commit message: "Changed method names and updated load test data method")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: rework the "Decision Process" test-scenario block.
# ---------------------------------------------------------------------------

# B3 / C3 become the start of a vertically-merged block (B3:B5 and C3:C5),
# centered both horizontally and vertically, still wrapping text.
$ws.Range("B3").Value = "Decision Process"
$ws.Range("C3").Value = "Test validates approval/rejection/review for loan"

# D3: method name rename.
$ws.Range("D3").Value = "TestAutoDecisionProcess"

# E3 is no longer used at all - remove value + formatting completely.
$ws.Range("E3").Clear()

# F3, G3, H3, I3 keep their existing values/formatting untouched.

# New alignment for the merged B/C block: center + vertical-center + wrap.
$ws.Range("B3:C5").HorizontalAlignment = -4108
$ws.Range("B3:C5").VerticalAlignment = -4108
$ws.Range("B3:C5").WrapText = $true

# ---------------------------------------------------------------------------
# Row 4 / Row 5: additional method names under the same scenario.
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "TestDisbursement"
$ws.Range("D5").Value = "TestManualReject"

# ---------------------------------------------------------------------------
# Merge the B/C columns across rows 3-5.
# ---------------------------------------------------------------------------
$ws.Range("B3:B5").Merge()
$ws.Range("C3:C5").Merge()

# ---------------------------------------------------------------------------
# Row 7 (row 6 intentionally left blank): new "Application" scenario.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = "Application"
$ws.Range("C7").Value = "Test validates addition and removal of applicants from system"
$ws.Range("D7").Value = "TestAddVerifyRemoveApplicant"

# ---------------------------------------------------------------------------
# Row heights for the wrapped, multi-line rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 45

# ---------------------------------------------------------------------------
# Column C is now wider to fit the longer descriptions.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 24.7

# ---------------------------------------------------------------------------
# Selection moves to F3 (matches the saved sheet view).
# ---------------------------------------------------------------------------
$ws.Range("F3").Select()

Write-Output "edit applied"
